$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is a valid numeric literal need to be forced
# to Text format first, otherwise Excel auto-converts the entered string into
# a genuine number (losing the literal "27.xxx"-style text representation that
# the source data/original workbook uses for the Price column).
$ws.Range("D2").Value = "27.296.70"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.832.11"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.21"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3681"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07441"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8846"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "1.898.07"
$ws.Range("E12").Value = "  +4.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07313"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.423"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.75"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.555"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008788"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "27.579.25"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.282"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "2.113.51"
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.894"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.81"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.54"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.134"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.221"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.14"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08977"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7494"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.174"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.538"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.942"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.096"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05334"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01955"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.423"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.953"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.233"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5286"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.481"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.48"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.07"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.009"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06293"
$ws.Range("E51").Value = "  +0.04%  "
